# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 43 - Afganistan
$ws.Range("B43").Value = 24102
$ws.Range("C43").Value = 556
$ws.Range("D43").Value = 4201
$ws.Range("E43").Value = 19450
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 451

# Row 48 - Israel
$ws.Range("B48").Value = 18876
$ws.Range("C48").Value = 81
$ws.Range("D48").Value = 15319
$ws.Range("E48").Value = 3257

# Row 76 - Uzbekistan
$ws.Range("B76").Value = 4901
$ws.Range("C76").Value = 32
$ws.Range("E76").Value = 1124

# Row 86 - El Salvador
$ws.Range("D86").Value = 1603
$ws.Range("E86").Value = 1806
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 72

# Row 90 - Tailandia
$ws.Range("B90").Value = 3134
$ws.Range("C90").Value = 5
$ws.Range("E90").Value = 89
